$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting "Total Debt/Equity (X)" and
# "Price/BV (X)" one column to the right. The new column inherits the header
# formatting (bold, centered, bordered) from the preceding column automatically.
$ws.Range("C1").EntireColumn.Insert()

# Header for the newly inserted column C
$ws.Range("C1").Value = "Book Value [ExclRevalReserve]/Share (Rs.)"

$values = @{
    2  = "203.66"
    3  = "210.86"
    4  = "202.15"
    5  = "198.31"
    6  = "210.39"
    7  = "397.20"
    8  = "396.05"
    9  = "330.02"
    10 = "231.87"
    11 = "224.90"
    12 = "165.86"
    13 = "126.49"
    14 = "99.53"
    15 = "76.73"
    16 = "136.38"
    17 = "111.43"
    18 = "82.35"
    19 = "114.64"
    20 = "69.17"
    21 = "12.92"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Range("C$row")
    $cell.Value = "'" + $values[$row]
    $cell.ClearFormats()
}
